$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value2 = 789.375
$ws.Range("I4").Value2 = 187.85715
$ws.Range("K4").Value2 = 187.85715
$ws.Range("M4").Value2 = -73.85714999999999
$ws.Range("H11").Value2 = 195.85715
$ws.Range("I11").Value2 = 195.85715
$ws.Range("K11").Value2 = 195.85715
$ws.Range("M11").Value2 = -55.85714999999999
$ws.Range("H17").Value2 = 872.64
$ws.Range("J17").Value2 = 872.64
$ws.Range("L17").Value2 = 2617.92
$ws.Range("N17").Value2 = -2953.92
$ws.Range("H19").Value2 = 683.3077
$ws.Range("I19").Value2 = 626.8570999999999
$ws.Range("K19").Value2 = 626.8570999999999
$ws.Range("M19").Value2 = -451.8570999999999
$ws.Range("H32").Value2 = 2171.484
$ws.Range("I32").Value2 = 1867.125
$ws.Range("J32").Value2 = 2277.348
$ws.Range("K32").Value2 = 1867.125
$ws.Range("L32").Value2 = 2277.348
$ws.Range("M32").Value2 = -1541.125
$ws.Range("N32").Value2 = -2929.348
$ws.Range("H39").Value2 = 2127.3333
$ws.Range("I39").Value2 = 532.8570999999999
$ws.Range("J39").Value2 = 4359.6
$ws.Range("K39").Value2 = 1598.5713
$ws.Range("L39").Value2 = 13078.8
$ws.Range("M39").Value2 = -1302.5713
$ws.Range("N39").Value2 = -13670.8
$ws.Range("H51").Value2 = 12479.733
$ws.Range("I51").Value2 = 12299.833
$ws.Range("J51").Value2 = 12599.667
$ws.Range("K51").Value2 = 12299.833
$ws.Range("L51").Value2 = 12599.667
$ws.Range("M51").Value2 = -11815.833
$ws.Range("N51").Value2 = -13567.667
$ws.Range("H70").Value2 = 16624.875
$ws.Range("I70").Value2 = 9000
$ws.Range("J70").Value2 = 19166.5
$ws.Range("K70").Value2 = 27000
$ws.Range("L70").Value2 = 57499.5
$ws.Range("M70").Value2 = -26730
$ws.Range("N70").Value2 = -58039.5
$ws.Range("H73").Value2 = 16624.875
$ws.Range("I73").Value2 = 9000
$ws.Range("J73").Value2 = 19166.5
$ws.Range("K73").Value2 = 27000
$ws.Range("L73").Value2 = 57499.5
$ws.Range("M73").Value2 = -26064
$ws.Range("N73").Value2 = -59371.5
$ws.Range("H74").Value2 = 7163
$ws.Range("I74").Value2 = 5974.5713
$ws.Range("K74").Value2 = 5974.5713
$ws.Range("M74").Value2 = -5038.5713
$ws.Range("H77").Value2 = 7163
$ws.Range("I77").Value2 = 5974.5713
$ws.Range("K77").Value2 = 29872.8565
$ws.Range("M77").Value2 = -25192.8565
$ws.Range("H98").Value2 = 0
$ws.Range("I98").Value2 = 0
$ws.Range("K98").Value2 = 0
$ws.Range("M98").ClearContents()
$ws.Range("H100").Value2 = 5579.8
$ws.Range("J100").Value2 = 6249.875
$ws.Range("L100").Value2 = 6249.875
$ws.Range("N100").Value2 = -7331.875
$ws.Range("H111").Value2 = 3357.5667
$ws.Range("I111").Value2 = 2903.4666
$ws.Range("J111").Value2 = 3811.6667
$ws.Range("K111").Value2 = 8710.399800000001
$ws.Range("L111").Value2 = 11435.0001
$ws.Range("M111").Value2 = -5643.399800000001
$ws.Range("N111").Value2 = -17569.0001
$ws.Range("H113").Value2 = 7010.923
$ws.Range("I113").Value2 = 6070.4546
$ws.Range("J113").Value2 = 7700.6
$ws.Range("K113").Value2 = 6070.4546
$ws.Range("L113").Value2 = 7700.6
$ws.Range("M113").Value2 = -2816.4546
$ws.Range("N113").Value2 = -14208.6
$ws.Range("H116").Value2 = 14159.8125
$ws.Range("I116").Value2 = 13432.091
$ws.Range("J116").Value2 = 15760.8
$ws.Range("K116").Value2 = 13432.091
$ws.Range("L116").Value2 = 15760.8
$ws.Range("M116").Value2 = -9990.091
$ws.Range("N116").Value2 = -22644.8
$ws.Range("H122").Value2 = 0
$ws.Range("I122").Value2 = 0
$ws.Range("K122").Value2 = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value2 = 2759.3708
$ws.Range("I132").Value2 = 2553.1206
$ws.Range("K132").Value2 = 7659.361800000001
$ws.Range("M132").Value2 = -5129.361800000001
$ws.Range("H137").Value2 = 2407.75
$ws.Range("I137").Value2 = 1515.1111
$ws.Range("K137").Value2 = 4545.3333
$ws.Range("M137").Value2 = -1995.3333
$ws.Range("H138").Value2 = 3211.561
$ws.Range("J138").Value2 = 4505.6
$ws.Range("L138").Value2 = 13516.8
$ws.Range("N138").Value2 = -23796.8

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1216
$ws.Range("I2").Value2 = 1208.1177
$ws.Range("K2").Value2 = 1208.1177
$ws.Range("M2").Value2 = -1095.1177
$ws.Range("H32").Value2 = 3291.5151
$ws.Range("I32").Value2 = 3191.25
$ws.Range("K32").Value2 = 3191.25
$ws.Range("M32").Value2 = -2904.25
$ws.Range("H45").Value2 = 2097.3333
$ws.Range("I45").Value2 = 1396.25
$ws.Range("K45").Value2 = 1396.25
$ws.Range("M45").Value2 = -1019.25
$ws.Range("H61").Value2 = 8612.210999999999
$ws.Range("I61").Value2 = 6600
$ws.Range("J61").Value2 = 16158
$ws.Range("K61").Value2 = 6600
$ws.Range("L61").Value2 = 16158
$ws.Range("M61").Value2 = -6388
$ws.Range("N61").Value2 = -16582
$ws.Range("H116").Value2 = 1216
$ws.Range("I116").Value2 = 1208.1177
$ws.Range("K116").Value2 = 1208.1177
$ws.Range("M116").Value2 = 1085.8823
$ws.Range("H132").Value2 = 1818.2273
$ws.Range("I132").Value2 = 1712.579
$ws.Range("K132").Value2 = 5137.737
$ws.Range("M132").Value2 = -2607.737
$ws.Range("H136").Value2 = 8612.210999999999
$ws.Range("I136").Value2 = 6600
$ws.Range("J136").Value2 = 16158
$ws.Range("K136").Value2 = 19800
$ws.Range("L136").Value2 = 48474
$ws.Range("M136").Value2 = -17250
$ws.Range("N136").Value2 = -53574

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1216
$ws.Range("I3").Value2 = 1208.1177
$ws.Range("K3").Value2 = 1208.1177
$ws.Range("M3").Value2 = -1094.1177
$ws.Range("H107").Value2 = 2602.2727
$ws.Range("I107").Value2 = 2582.2856
$ws.Range("K107").Value2 = 2582.2856
$ws.Range("M107").Value2 = -662.2856000000002
$ws.Range("H134").Value2 = 12442.737
$ws.Range("I134").Value2 = 15387.5
$ws.Range("K134").Value2 = 46162.5
$ws.Range("M134").Value2 = -43627.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3111.5925
$ws.Range("I31").Value2 = 1374
$ws.Range("J31").Value2 = 3328.7917
$ws.Range("K31").Value2 = 1374
$ws.Range("L31").Value2 = 3328.7917
$ws.Range("M31").Value2 = -1079
$ws.Range("N31").Value2 = -3918.7917
$ws.Range("H34").Value2 = 3111.5925
$ws.Range("I34").Value2 = 1374
$ws.Range("J34").Value2 = 3328.7917
$ws.Range("K34").Value2 = 1374
$ws.Range("L34").Value2 = 3328.7917
$ws.Range("M34").Value2 = -1172
$ws.Range("N34").Value2 = -3732.7917
$ws.Range("H52").Value2 = 81300
$ws.Range("J52").Value2 = 81300
$ws.Range("L52").Value2 = 81300
$ws.Range("N52").Value2 = -81888
$ws.Range("H132").Value2 = 2347.762
$ws.Range("I132").Value2 = 2361.7778
$ws.Range("J132").Value2 = 2263.6667
$ws.Range("K132").Value2 = 7085.3334
$ws.Range("L132").Value2 = 6791.000100000001
$ws.Range("M132").Value2 = -4555.3334
$ws.Range("N132").Value2 = -11851.0001
$ws.Range("H134").Value2 = 5974.7954
$ws.Range("I134").Value2 = 5207.5454
$ws.Range("K134").Value2 = 15622.6362
$ws.Range("M134").Value2 = -13087.6362

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 4492834
$ws.Range("I131").Value2 = 9261322
$ws.Range("J131").Value2 = 3835111.5
$ws.Range("K131").Value2 = 27783966
$ws.Range("L131").Value2 = 11505334.5
$ws.Range("M131").Value2 = -27778926
$ws.Range("N131").Value2 = -11515414.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value2 = 14414.286
$ws.Range("I46").Value2 = 11180
$ws.Range("J46").Value2 = 22500
$ws.Range("K46").Value2 = 11180
$ws.Range("L46").Value2 = 22500
$ws.Range("M46").Value2 = -11024
$ws.Range("N46").Value2 = -22812
$ws.Range("H102").Value2 = 6336.1113
$ws.Range("I102").Value2 = 6004.5
$ws.Range("K102").Value2 = 6004.5
$ws.Range("M102").Value2 = -4382.5
$ws.Range("H122").Value2 = 2604.524
$ws.Range("I122").Value2 = 1599.4
$ws.Range("J122").Value2 = 2918.625
$ws.Range("K122").Value2 = 4798.200000000001
$ws.Range("L122").Value2 = 8755.875
$ws.Range("M122").Value2 = -2348.200000000001
$ws.Range("N122").Value2 = -13655.875
$ws.Range("H132").Value2 = 2596.4443
$ws.Range("I132").Value2 = 2308
$ws.Range("K132").Value2 = 6924
$ws.Range("M132").Value2 = -4394

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value2 = 2646.1875
$ws.Range("I68").Value2 = 2583.9
$ws.Range("J68").Value2 = 2750
$ws.Range("K68").Value2 = 2583.9
$ws.Range("L68").Value2 = 2750
$ws.Range("M68").Value2 = -1834.9
$ws.Range("N68").Value2 = -4248
$ws.Range("H71").Value2 = 2646.1875
$ws.Range("I71").Value2 = 2583.9
$ws.Range("J71").Value2 = 2750
$ws.Range("K71").Value2 = 12919.5
$ws.Range("L71").Value2 = 13750
$ws.Range("M71").Value2 = -9175.5
$ws.Range("N71").Value2 = -21238
$ws.Range("H122").Value2 = 5879.95
$ws.Range("I122").Value2 = 4737.467
$ws.Range("K122").Value2 = 14212.401
$ws.Range("M122").Value2 = -11762.401
$ws.Range("H132").Value2 = 3049.1538
$ws.Range("I132").Value2 = 2714
$ws.Range("K132").Value2 = 8142
$ws.Range("M132").Value2 = -5612

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 2960.2942
$ws.Range("I132").Value2 = 2945.0715
$ws.Range("K132").Value2 = 8835.2145
$ws.Range("M132").Value2 = -6305.2145
